$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.58"
$ws.Range("E2").Value = "'-3.81%"
$ws.Range("D3").Value = "'37.17"
$ws.Range("E3").Value = "'-6.73%"
$ws.Range("D4").Value = "'5.098"
$ws.Range("E4").Value = "'-0.83%"
$ws.Range("D5").Value = "'0.07722"
$ws.Range("E5").Value = "'-6.06%"
$ws.Range("D6").Value = "'4.384"
$ws.Range("E6").Value = "'1.12%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.213"
$ws.Range("E7").Value = "'-1.59%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.866"
$ws.Range("E8").Value = "'-9.52%"
$ws.Range("E9").Value = "'-1.03%"
$ws.Range("D10").Value = "'0.9188"
$ws.Range("E10").Value = "'-1.97%"
$ws.Range("D11").Value = "'0.1232"
$ws.Range("E11").Value = "'-10.42%"
$ws.Range("D12").Value = "'0.1878"
$ws.Range("E12").Value = "'-5.46%"
$ws.Range("D13").Value = "'0.08761"
$ws.Range("E13").Value = "'-4.28%"
$ws.Range("D14").Value = "'0.03431"
$ws.Range("E14").Value = "'-2.19%"
$ws.Range("D15").Value = "'0.09715"
$ws.Range("E15").Value = "'-0.76%"
$ws.Range("D16").Value = "'0.001372"
$ws.Range("E16").Value = "'-2.54%"
$ws.Range("D17").Value = "'0.006173"
$ws.Range("E17").Value = "'-1.01%"
$ws.Range("D18").Value = "'3.556"
$ws.Range("E18").Value = "'-3.76%"
$ws.Range("E19").Value = "'-3.53%"
$ws.Range("D20").Value = "'0.1283"
$ws.Range("E20").Value = "'-2.07%"
$ws.Range("D21").Value = "'5.033"
$ws.Range("E21").Value = "'1.75%"
$ws.Range("D22").Value = "'0.2501"
$ws.Range("E22").Value = "'2.13%"
$ws.Range("D23").Value = "'0.02114"
$ws.Range("E23").Value = "'5,189.64%"
$ws.Range("D24").Value = "'0.04334"
$ws.Range("E24").Value = "'-0.44%"
$ws.Range("D25").Value = "'0.001221"
$ws.Range("E25").Value = "'-0.64%"
$ws.Range("D26").Value = "'0.004465"
$ws.Range("E26").Value = "'-7.51%"
$ws.Range("D27").Value = "'0.0001357"
$ws.Range("E27").Value = "'4.39%"
$ws.Range("D39").Value = "'0.02208"
$ws.Range("E39").Value = "'-0.68%"
$ws.Range("D40").Value = "'0.04909"
$ws.Range("E40").Value = "'-5.87%"
$ws.Range("D41").Value = "'0.007630"
$ws.Range("E41").Value = "'-1.36%"
$ws.Range("D42").Value = "'0.009920"
$ws.Range("E42").Value = "'2.40%"
$ws.Range("D43").Value = "'0.1332"
$ws.Range("E43").Value = "'-5.34%"
$ws.Range("D44").Value = "'0.002004"
$ws.Range("E44").Value = "'-2.18%"
$ws.Range("D45").Value = "'0.008807"
$ws.Range("E45").Value = "'-8.75%"
$ws.Range("D46").Value = "'0.00006964"
$ws.Range("E46").Value = "'5.26%"
$ws.Range("D47").Value = "'0.00000000754"
$ws.Range("E47").Value = "'0.64%"
$ws.Range("D48").Value = "'0.003013"
$ws.Range("E48").Value = "'2.51%"
$ws.Range("D49").Value = "'0.001307"
$ws.Range("E49").Value = "'-22.57%"
$ws.Range("D50").Value = "'0.00002111"
$ws.Range("E50").Value = "'0.64%"
$ws.Range("D51").Value = "'0.0002011"
$ws.Range("E51").Value = "'0.64%"
